$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds free-form text (prices like "46.189.80" use dots as
# thousands separators and are not valid numbers), so force text format
# before writing values to avoid Excel auto-converting them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "46.278.19"
$ws.Range("E2").Value = "  -1.11%  "

$ws.Range("D3").Value = "2.476.85"
$ws.Range("E3").Value = "  +9.40%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "297.32"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").Value = "95.68"
$ws.Range("E6").Value = "  -3.27%  "

$ws.Range("E7").Value = "  +1.66%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  +1.89%  "

$ws.Range("D10").Value = "35.22"
$ws.Range("E10").Value = "  +1.06%  "

$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").Value = "7.17"
$ws.Range("E12").Value = "  +2.35%  "

$ws.Range("E13").Value = "  +1.91%  "

$ws.Range("D14").Value = "2.855.94"
$ws.Range("E14").Value = "  +9.43%  "

$ws.Range("D15").Value = "2.462.65"
$ws.Range("E15").Value = "  +8.58%  "

$ws.Range("D16").Value = "0.863"
$ws.Range("E16").Value = "  +8.73%  "

$ws.Range("D17").Value = "14.24"
$ws.Range("E17").Value = "  +4.91%  "

$ws.Range("D18").Value = "46.305.14"
$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("D19").Value = "12.76"
$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  -2.00%  "

$ws.Range("D21").Value = "6.34"
$ws.Range("E21").Value = "  +9.52%  "

$ws.Range("D22").Value = "67.76"
$ws.Range("E22").Value = "  +3.25%  "

$ws.Range("D23").Value = "247.16"
$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("E24").Value = "  +1.74%  "

$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  +6.32%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "39.58"
$ws.Range("E27").Value = "  -3.62%  "

$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").Value = "9.87"
$ws.Range("E29").Value = "  +3.95%  "

$ws.Range("D30").Value = "3.90"
$ws.Range("E30").Value = "  +17.28%  "

$ws.Range("D31").Value = "21.72"
$ws.Range("E31").Value = "  +8.87%  "

$ws.Range("D32").Value = "2.76"
$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("D33").Value = "5.61"
$ws.Range("E33").Value = "  +5.87%  "

$ws.Range("D34").Value = "147.70"
$ws.Range("E34").Value = "  +1.13%  "

$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  +24.79%  "

$ws.Range("D36").Value = "0.0775"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("E37").Value = "  +2.91%  "

$ws.Range("E38").Value = "  +1.22%  "

$ws.Range("D39").Value = "15.27"
$ws.Range("E39").Value = "  -1.63%  "

$ws.Range("D40").Value = "3.97"
$ws.Range("E40").Value = "  +4.46%  "

$ws.Range("E41").Value = "  +2.69%  "

$ws.Range("D42").Value = "3.28"
$ws.Range("E42").Value = "  +7.45%  "

$ws.Range("D43").Value = "2.000.76"
$ws.Range("E43").Value = "  +12.36%  "

$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").Value = "92.71"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").Value = "16.57"
$ws.Range("E46").Value = "  +33.91%  "

$ws.Range("D47").Value = "1.79"
$ws.Range("E47").Value = "  -4.73%  "

$ws.Range("D48").Value = "8.61"
$ws.Range("E48").Value = "  +10.21%  "

$ws.Range("D49").Value = "103.05"
$ws.Range("E49").Value = "  +9.63%  "

$ws.Range("D50").Value = "2.716.23"
$ws.Range("E50").Value = "  +9.26%  "

$ws.Range("D51").Value = "0.186"
$ws.Range("E51").Value = "  +2.03%  "
